$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows continue the daily series (serial dates) started in the sheet,
# each with nuovi pos. / somma mobile columns at 0, matching the existing
# look/format (same style as the row above) -> "aggiornamento fino a 13/03"
$dates = @(44326, 44327, 44328, 44329)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 252 + $i
    $prevRow = $row - 1

    # Clone formatting (style) from the row above onto the new row first
    $ws.Range("A$prevRow`:D$prevRow").Copy()
    $ws.Range("A$row`:D$row").PasteSpecial(-4122)

    $ws.Range("A$row").Value = $dates[$i]
    $ws.Range("B$row").Value = 0
    $ws.Range("C$row").Value = 0
    $ws.Range("D$row").Value = 0
}
